# Update cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextCell "D2" "26.655.51"
Set-TextCell "E2" "  +0.18%  "
Set-TextCell "D3" "1.643.66"
Set-TextCell "E3" "  +0.86%  "
Set-TextCell "E4" "  +0.11%  "
Set-TextCell "D5" "215.48"
Set-TextCell "E5" "  +1.12%  "
Set-TextCell "E6" "  +1.33%  "
Set-TextCell "E7" "  +0.14%  "
Set-TextCell "E8" "  +0.53%  "
Set-TextCell "E9" "  +0.94%  "
Set-TextCell "D10" "19.27"
Set-TextCell "E10" "  +0.45%  "
Set-TextCell "E11" "  -0.19%  "
Set-TextCell "D12" "1.872.38"
Set-TextCell "E12" "  +0.77%  "
Set-TextCell "D13" "1.638.16"
Set-TextCell "E13" "  +0.21%  "
Set-TextCell "E14" "  +2.21%  "
Set-TextCell "E15" "  +1.83%  "
Set-TextCell "E16" "  +3.05%  "
Set-TextCell "D17" "26.696.28"
Set-TextCell "E17" "  +0.37%  "
Set-TextCell "E18" "  +0.55%  "
Set-TextCell "D19" "216.85"
Set-TextCell "E19" "  +0.50%  "
Set-TextCell "E20" "  +0.08%  "
Set-TextCell "D21" "4.35"
Set-TextCell "E21" "  +0.89%  "
Set-TextCell "E22" "  +2.49%  "
Set-TextCell "E23" "  +1.83%  "
Set-TextCell "D24" "2.26"
Set-TextCell "E24" "  +14.76%  "
Set-TextCell "D25" "145.69"
Set-TextCell "E25" "  -1.19%  "
Set-TextCell "E27" "  -0.45%  "
Set-TextCell "D28" "7.18"
Set-TextCell "E28" "  +4.91%  "
Set-TextCell "D29" "15.78"
Set-TextCell "E29" "  +1.80%  "
Set-TextCell "D30" "0.0516"
Set-TextCell "E30" "  +2.56%  "
Set-TextCell "E31" "  +0.78%  "
Set-TextCell "E32" "  +2.85%  "
Set-TextCell "E33" "  +2.53%  "
Set-TextCell "D34" "1.276.24"
Set-TextCell "E34" "  +4.23%  "
Set-TextCell "E35" "  +3.18%  "
Set-TextCell "E36" "  +5.05%  "
Set-TextCell "E37" "  +0.74%  "
Set-TextCell "D38" "0.533"
Set-TextCell "E38" "  +6.78%  "
Set-TextCell "D39" "0.829"
Set-TextCell "E39" "  +3.30%  "
Set-TextCell "E40" "  +0.25%  "
Set-TextCell "E41" "  +2.64%  "
Set-TextCell "E42" "  -1.61%  "
Set-TextCell "D43" "5.47"
Set-TextCell "E43" "  +2.45%  "
Set-TextCell "D44" "1.782.59"
Set-TextCell "E44" "  +0.86%  "
Set-TextCell "D45" "91.77"
Set-TextCell "E45" "  -1.11%  "
Set-TextCell "D46" "59.87"
Set-TextCell "E46" "  +8.70%  "
Set-TextCell "E47" "  +1.78%  "
Set-TextCell "E48" "  +0.77%  "
Set-TextCell "D49" "7.76"
Set-TextCell "E49" "  +2.09%  "
Set-TextCell "E50" "  +3.41%  "
Set-TextCell "D51" "0.406"
Set-TextCell "E51" "  -0.82%  "
